$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# "hieight fix": the eight 86.3pt (1726-twip) data rows in the first
# repeated block shrink to 75.1pt (1502 twips). These are rows 2-9 of
# the single big table (row 1 is the header, row 10 begins the next
# section) - row 17 (the very last row) is also 86.3pt but is left
# untouched, matching the diff.
for ($i = 2; $i -le 9; $i++) {
    $t.Rows.Item($i).Height = 75.1
}

# The last 43.15pt (863-twip) row (row 16, the "Notes" row right before
# the final 86.3pt row) shrinks to 42.75pt (855 twips).
$t.Rows.Item(16).Height = 42.75
